{"js": "const pairs = [\n  [\"61\u00d789=5429\", \"73\u00d732=2336\"],\n  [\"36\u00d753=1908\", \"66\u00d744=2904\"],\n  [\"32\u00d716=512\", \"84\u00d757=4788\"],\n  [\"80\u00d772=5760\", \"80\u00d761=4880\"],\n  [\"79\u00d723=1817\", \"52\u00d771=3692\"],\n  [\"76\u00d779=6004\", \"35\u00d717=595\"],\n  [\"62\u00d771=4402\", \"91\u00d716=1456\"],\n  [\"31\u00d731=961\", \"42\u00d711=462\"],\n  [\"71\u00d792=6532\", \"14\u00d799=1386\"],\n  [\"81\u00d762=5022\", \"55\u00d712=660\"],\n  [\"97\u00d749=4753\", \"77\u00d711=847\"],\n  [\"82\u00d716=1312\", \"63\u00d762=3906\"],\n  [\"15\u00d794=1410\", \"64\u00d761=3904\"],\n  [\"25\u00d721=525\", \"65\u00d768=4420\"],\n  [\"34\u00d754=1836\", \"66\u00d757=3762\"],\n  [\"33\u00d762=2046\", \"94\u00d726=2444\"],\n  [\"19\u00d747=893\", \"55\u00d723=1265\"],\n  [\"13\u00d782=1066\", \"31\u00d737=1147\"],\n  [\"78\u00d754=4212\", \"78\u00d711=858\"],\n  [\"94\u00d775=7050\", \"88\u00d788=7744\"],\n  [\"78\u00d751=3978\", \"19\u00d737=703\"],\n  [\"93\u00d790=8370\", \"91\u00d726=2366\"],\n  [\"66\u00d761=4026\", \"85\u00d751=4335\"],\n  [\"13\u00d776=988\", \"89\u00d770=6230\"],\n  [\"13\u00d725=325\", \"40\u00d761=2440\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWildcards: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n  @(\"61\u00d789=5429\", \"73\u00d732=2336\"),\n  @(\"36\u00d753=1908\", \"66\u00d744=2904\"),\n  @(\"32\u00d716=512\", \"84\u00d757=4788\"),\n  @(\"80\u00d772=5760\", \"80\u00d761=4880\"),\n  @(\"79\u00d723=1817\", \"52\u00d771=3692\"),\n  @(\"76\u00d779=6004\", \"35\u00d717=595\"),\n  @(\"62\u00d771=4402\", \"91\u00d716=1456\"),\n  @(\"31\u00d731=961\", \"42\u00d711=462\"),\n  @(\"71\u00d792=6532\", \"14\u00d799=1386\"),\n  @(\"81\u00d762=5022\", \"55\u00d712=660\"),\n  @(\"97\u00d749=4753\", \"77\u00d711=847\"),\n  @(\"82\u00d716=1312\", \"63\u00d762=3906\"),\n  @(\"15\u00d794=1410\", \"64\u00d761=3904\"),\n  @(\"25\u00d721=525\", \"65\u00d768=4420\"),\n  @(\"34\u00d754=1836\", \"66\u00d757=3762\"),\n  @(\"33\u00d762=2046\", \"94\u00d726=2444\"),\n  @(\"19\u00d747=893\", \"55\u00d723=1265\"),\n  @(\"13\u00d782=1066\", \"31\u00d737=1147\"),\n  @(\"78\u00d754=4212\", \"78\u00d711=858\"),\n  @(\"94\u00d775=7050\", \"88\u00d788=7744\"),\n  @(\"78\u00d751=3978\", \"19\u00d737=703\"),\n  @(\"93\u00d790=8370\", \"91\u00d726=2366\"),\n  @(\"66\u00d761=4026\", \"85\u00d751=4335\"),\n  @(\"13\u00d776=988\", \"89\u00d770=6230\"),\n  @(\"13\u00d725=325\", \"40\u00d761=2440\"),\n)\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
